# Fruta / hortaliza, semanal
# Insert a new weekly record at row 24 (pushing the existing rows 24-52
# down to 25-53) in the "Agrícola del Norte S.A. de Arica - Guayaba" log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 24, shifting everything
# below it (rows 24-52) down by one (to 25-53).
$ws.Rows("24:24").Insert()

# Fill in the new weekly record.
$ws.Range("A24").Value = 1
$ws.Range("B24").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C24").Value = "Arica y Parinacota"
$ws.Range("D24").Value = 44883
$ws.Range("D24").NumberFormat = $ws.Range("D25").NumberFormat
$ws.Range("E24").Value = 15
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100108
$ws.Range("H24").Value = "Tropicales y subtropicales"
$ws.Range("I24").Value = 100108001
$ws.Range("J24").Value = "Guayaba"
$ws.Range("K24").Value = "Sin especificar"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 100
$ws.Range("N24").Value = 700
$ws.Range("O24").Value = 750
$ws.Range("P24").Value = 725
$ws.Range("Q24").Value = "$/kilo (en caja de 10 kilos )"
$ws.Range("R24").Value = "Región de Arica y Parinacota"
$ws.Range("S24").Value = 725
$ws.Range("T24").Value = 1
